$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 43.859649
$ws.Range("D2").Value = 0.363008
$ws.Range("E2").Value = 0.547252

$ws.Range("B3").Value = 1598.171065
$ws.Range("D3").Value = 6.61369
$ws.Range("E3").Value = 0.001526

$ws.Range("B4").Value = 40113.22197
$ws.Range("C4").Value = 332

$ws.Range("G5").Value = 3.742149
$ws.Range("H5").Value = -0.243131
$ws.Range("I5").Value = 7.727428
$ws.Range("J5").Value = 0.07084500000000001

$ws.Range("G6").Value = -0.877193
$ws.Range("H6").Value = -5.063088
$ws.Range("I6").Value = 3.308702
$ws.Range("J6").Value = 0.8745270000000001

$ws.Range("G7").Value = -4.619342
$ws.Range("H7").Value = -7.766855
$ws.Range("I7").Value = -1.471828
$ws.Range("J7").Value = 0.001794
